# Update generated output (commit 456a3b4): bump "want-to-go" counts (column F)
# across sheets, and insert a new event row ("今泉爱夏 巡演") into the
# "演出" (Performances) and "全部类型" (All types) sheets.

$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("展览")
$ws2 = $wb.Worksheets.Item("演出")
$ws4 = $wb.Worksheets.Item("全部类型")

# ---------------------------------------------------------------
# Sheet "展览" (Exhibitions) - column F value bumps
# ---------------------------------------------------------------
$ws1.Range("F2").Value = 101
$ws1.Range("F4").Value = 412
$ws1.Range("F5").Value = 186
$ws1.Range("F6").Value = 129
$ws1.Range("F7").Value = 1107
$ws1.Range("F8").Value = 372
$ws1.Range("F10").Value = 49
$ws1.Range("F12").Value = 371
$ws1.Range("F13").Value = 377
$ws1.Range("F14").Value = 779
$ws1.Range("F15").Value = 162
$ws1.Range("F16").Value = 718
$ws1.Range("F17").Value = 273
$ws1.Range("F19").Value = 992
$ws1.Range("F20").Value = 448
$ws1.Range("F21").Value = 258
$ws1.Range("F22").Value = 81
$ws1.Range("F23").Value = 375
$ws1.Range("F24").Value = 26
$ws1.Range("F25").Value = 39
$ws1.Range("F26").Value = 464

# ---------------------------------------------------------------
# Sheet "演出" (Performances) - column F value bumps
# ---------------------------------------------------------------
$ws2.Range("F4").Value = 361
$ws2.Range("F5").Value = 39
$ws2.Range("F10").Value = 627
$ws2.Range("F11").Value = 148
$ws2.Range("F12").Value = 33

# Insert a new row 13 (shifts the old row 13 "夏川里美" entry down to row 14)
$ws2.Rows.Item(13).Insert()

# The row that shifted down keeps its old index value in column A; bump it by
# one to preserve the running sequence number (12 -> 13).
$ws2.Range("A14").Value = 13

# Copy formatting for the new row 13 from row 12 (an existing, fully
# formatted data row) so the styling (e.g. bold index column) matches the
# rest of the table.
$ws2.Range("A12:I12").Copy()
$ws2.Range("A13:I13").PasteSpecial(-4122)
$ws2.Range("A13").Value = 12

# Column B holds a plain date-like label ("2024.04.24"); force text storage
# so Excel doesn't reinterpret it as a date serial number, then reapply the
# plain (un-formatted) style used by the rest of the table.
$ws2.Range("B13").NumberFormat = "@"
$ws2.Range("B13").Value = "2024.04.24"
$ws2.Range("B12").Copy()
$ws2.Range("B13").PasteSpecial(-4122)

$ws2.Range("C13").Value = "广州·今泉爱夏  巡演"
$ws2.Range("D13").Value = "革新路124号太古仓码头54汇5号仓 太空间Livehouse"
$ws2.Range("E13").Value = "2024.04.24 20:00-04.24 21:30"
$ws2.Range("F13").Value = 0
$ws2.Range("G13").Value = "不可售"
$ws2.Range("H13").Value = "https://show.bilibili.com/platform/detail.html?id=81890"
$ws2.Range("I13").Value = "//i1.hdslb.com/bfs/openplatform/202402/YJENeaUi1708313389899.jpeg"

# ---------------------------------------------------------------
# Sheet "全部类型" (All types) - column F value bumps
# ---------------------------------------------------------------
$ws4.Range("F4").Value = 101
$ws4.Range("F6").Value = 412
$ws4.Range("F7").Value = 186
$ws4.Range("F8").Value = 129
$ws4.Range("F9").Value = 1107
$ws4.Range("F10").Value = 372
$ws4.Range("F13").Value = 49
$ws4.Range("F14").Value = 361
$ws4.Range("F16").Value = 39
$ws4.Range("F17").Value = 371
$ws4.Range("F20").Value = 377
$ws4.Range("F21").Value = 779
$ws4.Range("F22").Value = 162
$ws4.Range("F23").Value = 718
$ws4.Range("F24").Value = 273
$ws4.Range("F26").Value = 992
$ws4.Range("F27").Value = 448
$ws4.Range("F30").Value = 258
$ws4.Range("F31").Value = 81
$ws4.Range("F32").Value = 375
$ws4.Range("F34").Value = 148
$ws4.Range("F35").Value = 26
$ws4.Range("F36").Value = 39
$ws4.Range("F37").Value = 33
$ws4.Range("F38").Value = 464

# Insert a new row 39 (shifts the old row 39 "夏川里美" entry down to row 40)
$ws4.Rows.Item(39).Insert()

# The row that shifted down keeps its old index value in column A; bump it by
# one to preserve the running sequence number (38 -> 39).
$ws4.Range("A40").Value = 39

$ws4.Range("A38:I38").Copy()
$ws4.Range("A39:I39").PasteSpecial(-4122)
$ws4.Range("A39").Value = 38

$ws4.Range("B39").NumberFormat = "@"
$ws4.Range("B39").Value = "2024.04.24"
$ws4.Range("B38").Copy()
$ws4.Range("B39").PasteSpecial(-4122)

$ws4.Range("C39").Value = "广州·今泉爱夏  巡演"
$ws4.Range("D39").Value = "革新路124号太古仓码头54汇5号仓 太空间Livehouse"
$ws4.Range("E39").Value = "2024.04.24 20:00-04.24 21:30"
$ws4.Range("F39").Value = 0
$ws4.Range("G39").Value = "不可售"
$ws4.Range("H39").Value = "https://show.bilibili.com/platform/detail.html?id=81890"
$ws4.Range("I39").Value = "//i1.hdslb.com/bfs/openplatform/202402/YJENeaUi1708313389899.jpeg"
